# Update "想去人数" (want-to-go count) values on sheets "展览", "演出", and "全部类型"
# to match the newly scraped data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1827
$ws1.Range("F6").Value  = 1436
$ws1.Range("F8").Value  = 1708
$ws1.Range("F22").Value = 4397
$ws1.Range("F23").Value = 31
$ws1.Range("F26").Value = 2117
$ws1.Range("F28").Value = 1995

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 69
$ws2.Range("F3").Value = 2

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1827
$ws4.Range("F6").Value  = 1436
$ws4.Range("F8").Value  = 1708
$ws4.Range("F22").Value = 4397
$ws4.Range("F23").Value = 69
$ws4.Range("F24").Value = 31
$ws4.Range("F25").Value = 2
$ws4.Range("F28").Value = 2117
$ws4.Range("F30").Value = 1995
